# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Uganda / Sri Lanka ordering (rows 127-128) ---
# Row 127 becomes "Uganda" (with Uganda's updated stats)
# Row 128 becomes "Sri Lanka" (keeping Sri Lanka's stats)
$ws.Range("A127").Value = "Uganda"
$ws.Range("A128").Value = "Sri Lanka"

# --- Update the "Datos actualizados" timestamp footer ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 14:21"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 6258805
$ws.Range("C4").Value = 1234
$ws.Range("D4").Value = 3497481
$ws.Range("E4").Value = 2572409
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 188915

# --- Row 64: Nepal ---
$ws.Range("B64").Value = 41649
$ws.Range("C64").Value = 1120
$ws.Range("D64").Value = 23290
$ws.Range("E64").Value = 18108
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 12
$ws.Range("H64").Value = 251

# --- Row 75: Estado de Palestina ---
$ws.Range("B75").Value = 23875
$ws.Range("C75").Value = 594
$ws.Range("D75").Value = 15483
$ws.Range("E75").Value = 8230
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 162

# --- Row 82: Dinamarca ---
$ws.Range("B82").Value = 17195
$ws.Range("C82").Value = 111
$ws.Range("D82").Value = 15413
$ws.Range("E82").Value = 1156
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 626

# --- Row 85: Madagascar ---
$ws.Range("B85").Value = 15023
$ws.Range("C85").Value = 66
$ws.Range("D85").Value = 13965
$ws.Range("E85").Value = 862
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 196

# --- Row 127: now Uganda's data ---
$ws.Range("B127").Value = 3112
$ws.Range("C127").Value = 75
$ws.Range("D127").Value = 1528
$ws.Range("E127").Value = 1552
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 32

# --- Row 128: now Sri Lanka's data ---
$ws.Range("B128").Value = 3092
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 2883
$ws.Range("E128").Value = 197
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 12

# --- Row 141: Islandia ---
$ws.Range("B141").Value = 2121
$ws.Range("C141").Value = 5
$ws.Range("D141").Value = 2016
$ws.Range("E141").Value = 95
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 10
